$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A's stored width goes from 9.7109375 -> 11.7109375 (now matching B/C).
# Excel's ColumnWidth setter snaps to an integer pixel grid (steps of 1/6 in
# the stored "character" width units), so 11.7109375 itself is not a
# reachable grid point; 10.83 is the closest input that lands the stored
# width on the nearest attainable grid value (11.666666666666666).
$ws.Columns.Item(1).ColumnWidth = 10.83

$ws.Range("A1").Value = 148.89827879314171
$ws.Range("B1").Value = 4.4757474183981794
$ws.Range("C1").Value = 1.1219178082191781
